$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 338, shifting existing rows 338..366 down to 339..367
$ws.Rows.Item(338).Insert()

# Populate the newly inserted row 338 with the new record
$ws.Cells.Item(338, 1).Value = 8
$ws.Cells.Item(338, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(338, 3).Value = "Coquimbo"
$ws.Cells.Item(338, 4).Value = 44769
$ws.Cells.Item(338, 4).NumberFormat = $ws.Cells.Item(339, 4).NumberFormat
$ws.Cells.Item(338, 5).Value = 4
$ws.Cells.Item(338, 6).Value = 100112032
$ws.Cells.Item(338, 7).Value = "Zapallo italiano"
$ws.Cells.Item(338, 8).Value = "Sin especificar"
$ws.Cells.Item(338, 9).Value = "Primera"
$ws.Cells.Item(338, 10).Value = 400
$ws.Cells.Item(338, 11).Value = 16500
$ws.Cells.Item(338, 12).Value = 17000
$ws.Cells.Item(338, 13).Value = 16750
$ws.Cells.Item(338, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(338, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(338, 16).Value = 335
$ws.Cells.Item(338, 17).Value = 50
$ws.Cells.Item(338, 18).Value = "Hortaliza"
